$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Append 9 new user_detail_h rows (22..30) below the existing 21 rows.
# We insert by copying the last existing row (21) downward so the new
# rows inherit the exact same per-column cell styles already used in
# the sheet (email column uses style index "2", is_active uses "1").
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 9; $i++) {
  $ws.Rows("21:21").Copy()
  $ws.Rows("22:22").Insert(-4121)   # xlShiftDown
}
$excel.CutCopyMode = $false

# Columns: A=id, B=uin, C=name, D=email, E=mobile, F=status_code,
# G=lang_code, H=last_login_method, I=is_active, J=cr_by,
# K=cr_dtimes, L=eff_dtimes
$rows = @(
  @{ r=22; a=110021; b=7316931025; c="Magdalena Weber";   d="magdalena.weber@xyz.com";   e=932122450 },
  @{ r=23; a=110022; b=9137847236; c="Adrienne Hoffman";  d="adrienne.hoffman@xyz.com";  e=848488000 },
  @{ r=24; a=110023; b=8428758532; c="Adrienne Mcgee";    d="adrienne.mcgee@xyz.com";    e=894773246 },
  @{ r=25; a=110024; b=9804209494; c="Amare Coleman";     d="amare.coleman@xyz.com";     e=956554588 },
  @{ r=26; a=110025; b=7105248214; c="Dawson Ibarra";     d="dawson.ibarra@xyz.com";     e=765455583 },
  @{ r=27; a=110026; b=9316557128; c="Elvis Mcmillan";    d="elvis.mcmillan@xyz.com";    e=884282274 },
  @{ r=28; a=110027; b=8103486949; c="Steve George";      d="steve.george@xyz.com";      e=971073663 },
  @{ r=29; a=110028; b=9601932866; c="Colton Elliott";    d="colton.elliott@xyz.com";    e=809908673 },
  @{ r=30; a=110029; b=9317596765; c="Carolyn Rodriguez"; d="carolyn.rodriguez@xyz.com"; e=818876429 }
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r, 1).Value = $row.a       # A - id
  $ws.Cells.Item($r, 2).Value = $row.b       # B - uin
  $ws.Cells.Item($r, 3).Value = $row.c       # C - name
  $ws.Cells.Item($r, 4).Value = $row.d       # D - email
  $ws.Cells.Item($r, 5).Value = $row.e       # E - mobile
  $ws.Cells.Item($r, 6).Value = "ACT"        # F - status_code
  $ws.Cells.Item($r, 7).Value = "eng"        # G - lang_code
  $ws.Cells.Item($r, 8).Value = "PWD"        # H - last_login_method
  $ws.Cells.Item($r, 9).Value = $true        # I - is_active
  $ws.Cells.Item($r, 10).Value = "superadmin" # J - cr_by
  $ws.Cells.Item($r, 11).Value = "now()"      # K - cr_dtimes
  $ws.Cells.Item($r, 12).Value = "now()"      # L - eff_dtimes
}

# Move the visible viewport/selection to the newly appended rows.
$ws.Range("A22:A30").Select()
